$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the E1 (Category *) comment text to describe the new "Instrument" column
$null = $ws.Range("E1").Comment.Text("Author:`n-Mandatory`n-This is the instrument of the portfolio company being valued")

# 2. Remove the comment that lived on F1 (Sub Category *) - that column is going away
$ws.Range("F1").Comment.Delete() | Out-Null

# 3. Repurpose column E from "Category *"/"Unlisted" to "Instrument"/"Common Stock"
$ws.Range("E1").Value = "Instrument"
$ws.Range("E2:E5").Value = "Common Stock"

# 4. Drop column F ("Sub Category *" / "Equity") entirely - no longer needed
$ws.Range("F1").EntireColumn.Delete() | Out-Null

# 5. Move the active selection
$ws.Range("G10").Select() | Out-Null
